$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.164.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.96%  "
$ws.Range("D3").Value = "'3.570.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.74%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'416.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "'129.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").Value = "'0.651"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.67%  "
$ws.Range("D8").Value = "'3.560.28"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.68%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "'0.771"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.94%  "
$ws.Range("E11").Value = "  +14.19%  "
$ws.Range("D12").Value = "'0.0000336"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +49.05%  "
$ws.Range("D13").Value = "'42.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").Value = "'10.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.63%  "
$ws.Range("D15").Value = "'4.137.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.66%  "
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "'20.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "'3.571.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.16%  "
$ws.Range("E19").Value = "  +5.64%  "
$ws.Range("D20").Value = "'67.032.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.80%  "
$ws.Range("D21").Value = "'12.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.36%  "
$ws.Range("D22").Value = "'456.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("D23").Value = "'88.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.60%  "
$ws.Range("E24").Value = "  -5.74%  "
$ws.Range("D25").Value = "'13.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("D26").Value = "'3.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.66%  "
$ws.Range("D27").Value = "'10.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.96%  "
$ws.Range("D28").Value = "'34.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.25%  "
$ws.Range("D29").Value = "'4.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.79%  "
$ws.Range("E30").Value = "  +4.26%  "
$ws.Range("D31").Value = "'12.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.27%  "
$ws.Range("E32").Value = "  +4.89%  "
$ws.Range("D33").Value = "'7.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.28%  "
$ws.Range("D34").Value = "'0.162"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.12%  "
$ws.Range("D35").Value = "'41.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").Value = "'56.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").Value = "'0.0₃0724"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +23.04%  "
$ws.Range("E40").Value = "  +8.93%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").Value = "'148.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("E45").Value = "  -2.92%  "
$ws.Range("E46").Value = "  -3.23%  "
$ws.Range("D47").Value = "'4.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.20%  "
$ws.Range("E48").Value = "  -3.76%  "
$ws.Range("E49").Value = "  -1.46%  "
$ws.Range("D50").Value = "'2.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +15.04%  "
$ws.Range("D51").Value = "'15.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.83%  "
